$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content edits: rename two stat rows
$ws.Range("A5").Value = "Sleep Resist"
$ws.Range("A4").Value = "Attack"

# New column D formatting: center alignment on populated D cells
$ws.Range("D1:D6").HorizontalAlignment = -4108
$ws.Range("D8:D9").HorizontalAlignment = -4108

# Column D width (new column, roughly 32.25 chars)
$ws.Columns(4).ColumnWidth = 31.5

# View changes: zoom + selection
$excel.ActiveWindow.Zoom = 85
$ws.Range("G11").Select()
